$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 22
$ws.Range("I5").Value = 19.333334
$ws.Range("J5").Value = 30
$ws.Range("K5").Value = 19.333334
$ws.Range("L5").Value = 30
$ws.Range("M5").Value = 95.66666599999999
$ws.Range("N5").Value = -260
$ws.Range("H26").Value = 11666.667
$ws.Range("I26").Value = 10000
$ws.Range("K26").Value = 10000
$ws.Range("M26").Value = -9656
$ws.Range("H58").Value = 323579.47
$ws.Range("I58").Value = 606341.4
$ws.Range("J58").Value = 2259.0908
$ws.Range("K58").Value = 1819024.2
$ws.Range("L58").Value = 6777.2724
$ws.Range("M58").Value = -1818874.2
$ws.Range("N58").Value = -7077.2724
$ws.Range("H62").Value = 1900
$ws.Range("I62").Value = 1900
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1900
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1276
$ws.Range("N62").ClearContents()
$ws.Range("H64").Value = 86541.414
$ws.Range("I64").Value = 202639.6
$ws.Range("J64").Value = 3614.1428
$ws.Range("K64").Value = 202639.6
$ws.Range("L64").Value = 3614.1428
$ws.Range("M64").Value = -202391.6
$ws.Range("N64").Value = -4110.1428
$ws.Range("H65").Value = 1900
$ws.Range("I65").Value = 1900
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 9500
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -6380
$ws.Range("N65").ClearContents()
$ws.Range("H67").Value = 86541.414
$ws.Range("I67").Value = 202639.6
$ws.Range("J67").Value = 3614.1428
$ws.Range("K67").Value = 202639.6
$ws.Range("L67").Value = 3614.1428
$ws.Range("M67").Value = -201781.6
$ws.Range("N67").Value = -5330.1428
$ws.Range("H70").Value = 1186.1818
$ws.Range("I70").Value = 999.2
$ws.Range("J70").Value = 1342
$ws.Range("K70").Value = 2997.6
$ws.Range("L70").Value = 4026
$ws.Range("M70").Value = -2727.6
$ws.Range("N70").Value = -4566
$ws.Range("H73").Value = 1186.1818
$ws.Range("I73").Value = 999.2
$ws.Range("J73").Value = 1342
$ws.Range("K73").Value = 2997.6
$ws.Range("L73").Value = 4026
$ws.Range("M73").Value = -2061.6
$ws.Range("N73").Value = -5898
$ws.Range("H76").Value = 3342.8572
$ws.Range("I76").Value = 3342.8572
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3342.8572
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3027.8572
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 3342.8572
$ws.Range("I79").Value = 3342.8572
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3342.8572
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2250.8572
$ws.Range("N79").ClearContents()
$ws.Range("H82").Value = 3007.1538
$ws.Range("J82").Value = 3999.7778
$ws.Range("L82").Value = 11999.3334
$ws.Range("N82").Value = -12811.3334
$ws.Range("H85").Value = 3007.1538
$ws.Range("J85").Value = 3999.7778
$ws.Range("L85").Value = 11999.3334
$ws.Range("N85").Value = -14807.3334
$ws.Range("H88").Value = 4917.222
$ws.Range("I88").Value = 1000
$ws.Range("J88").Value = 6036.4287
$ws.Range("K88").Value = 1000
$ws.Range("L88").Value = 6036.4287
$ws.Range("M88").Value = -594
$ws.Range("N88").Value = -6848.4287
$ws.Range("H91").Value = 4917.222
$ws.Range("I91").Value = 1000
$ws.Range("J91").Value = 6036.4287
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 6036.4287
$ws.Range("M91").Value = 404
$ws.Range("N91").Value = -8844.4287
$ws.Range("H97").Value = 500050
$ws.Range("I97").Value = 150
$ws.Range("J97").Value = 750000
$ws.Range("K97").Value = 450
$ws.Range("L97").Value = 2250000
$ws.Range("M97").Value = 46
$ws.Range("N97").Value = -2250992
$ws.Range("H100").Value = 1701.5
$ws.Range("I100").Value = 602
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 602
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -61
$ws.Range("N100").Value = -6082
$ws.Range("H103").Value = 905.25
$ws.Range("I103").Value = 590
$ws.Range("J103").Value = 912.5814
$ws.Range("K103").Value = 1770
$ws.Range("L103").Value = 2737.7442
$ws.Range("M103").Value = -1184
$ws.Range("N103").Value = -3909.7442
$ws.Range("H106").Value = 2901.3333
$ws.Range("I106").Value = 2901.3333
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 2901.3333
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -2270.3333
$ws.Range("N106").ClearContents()
$ws.Range("H132").Value = 6584708.5
$ws.Range("I132").Value = 6762606
$ws.Range("K132").Value = 20287818
$ws.Range("M132").Value = -20285288

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H74").Value = 2245.7805
$ws.Range("I74").Value = 1634.1923
$ws.Range("J74").Value = 3305.8667
$ws.Range("K74").Value = 1634.1923
$ws.Range("L74").Value = 3305.8667
$ws.Range("M74").Value = -760.1922999999999
$ws.Range("N74").Value = -5053.8667
$ws.Range("H77").Value = 2245.7805
$ws.Range("I77").Value = 1634.1923
$ws.Range("J77").Value = 3305.8667
$ws.Range("K77").Value = 8170.961499999999
$ws.Range("L77").Value = 16529.3335
$ws.Range("M77").Value = -3802.961499999999
$ws.Range("N77").Value = -25265.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 406.5
$ws.Range("I10").Value = 406.5
$ws.Range("K10").Value = 406.5
$ws.Range("M10").Value = -267.5
$ws.Range("H122").Value = 4386.4287
$ws.Range("I122").Value = 3468
$ws.Range("J122").Value = 7325.4
$ws.Range("K122").Value = 10404
$ws.Range("L122").Value = 21976.2
$ws.Range("M122").Value = -7954
$ws.Range("N122").Value = -26876.2
$ws.Range("H124").Value = 39792
$ws.Range("J124").Value = 39792
$ws.Range("L124").Value = 39792
$ws.Range("N124").Value = -44702

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 1271629.2
$ws.Range("J37").Value = 1271629.2
$ws.Range("L37").Value = 3814887.6
$ws.Range("N37").Value = -3815111.6
$ws.Range("H112").Value = 1960.8334
$ws.Range("I112").Value = 950
$ws.Range("J112").Value = 2466.25
$ws.Range("K112").Value = 2850
$ws.Range("L112").Value = 7398.75
$ws.Range("M112").Value = -1742
$ws.Range("N112").Value = -9614.75
$ws.Range("H131").Value = 681329.9
$ws.Range("I131").Value = 602.8570999999999
$ws.Range("J131").Value = 733693.4399999999
$ws.Range("K131").Value = 1808.5713
$ws.Range("L131").Value = 2201080.32
$ws.Range("M131").Value = 3231.4287
$ws.Range("N131").Value = -2211160.32

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 4681818
$ws.Range("J7").Value = 3250000
$ws.Range("L7").Value = 3250000
$ws.Range("N7").Value = -3250224
$ws.Range("H8").Value = 4681818
$ws.Range("J8").Value = 3250000
$ws.Range("L8").Value = 3250000
$ws.Range("N8").Value = -3250278
$ws.Range("H14").Value = 2666736
$ws.Range("I14").Value = 2666736
$ws.Range("K14").Value = 2666736
$ws.Range("M14").Value = -2666568
$ws.Range("H47").Value = 4354
$ws.Range("J47").Value = 4354
$ws.Range("L47").Value = 4354
$ws.Range("N47").Value = -5490
$ws.Range("H134").Value = 26533.1
$ws.Range("J134").Value = 26533.1
$ws.Range("L134").Value = 79599.29999999999
$ws.Range("N134").Value = -84669.29999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 5143.3335
$ws.Range("I9").Value = 212
$ws.Range("J9").Value = 29800
$ws.Range("K9").Value = 212
$ws.Range("L9").Value = 29800
$ws.Range("M9").Value = 12
$ws.Range("N9").Value = -30248
$ws.Range("H20").Value = 42603.6
$ws.Range("J20").Value = 42603.6
$ws.Range("L20").Value = 42603.6
$ws.Range("N20").Value = -43055.6
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("H46").Value = 6499.75
$ws.Range("I46").Value = 1999
$ws.Range("J46").Value = 8000
$ws.Range("K46").Value = 1999
$ws.Range("L46").Value = 8000
$ws.Range("M46").Value = -1811
$ws.Range("N46").Value = -8376
$ws.Range("H55").Value = 944.56757
$ws.Range("I55").Value = 291.25
$ws.Range("J55").Value = 1124.7931
$ws.Range("K55").Value = 291.25
$ws.Range("L55").Value = 1124.7931
$ws.Range("M55").Value = -118.25
$ws.Range("N55").Value = -1470.7931
$ws.Range("H119").Value = 38485
$ws.Range("J119").Value = 38485
$ws.Range("L119").Value = 38485
$ws.Range("N119").Value = -48161

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 5221
$ws.Range("J54").Value = 4912.8335
$ws.Range("L54").Value = 4912.8335
$ws.Range("N54").Value = -5952.8335
$ws.Range("H109").Value = 37000
$ws.Range("J109").Value = 37000
$ws.Range("L109").Value = 37000
$ws.Range("N109").Value = -39774
$ws.Range("H124").Value = 13709.5
$ws.Range("J124").Value = 13709.5
$ws.Range("L124").Value = 13709.5
$ws.Range("N124").Value = -23529.5
$ws.Range("H140").Value = 51199.07
$ws.Range("J140").Value = 51199.07
$ws.Range("L140").Value = 51199.07
$ws.Range("N140").Value = -61559.07
